# Auto-generated edits applying the diff to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.176.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.263.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -7.37%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.54%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.255.94'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.82%  '
$ws.Range('E10').Value = '  -14.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.507'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -16.59%  '
$ws.Range('E14').Value = '  -11.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.784.99'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.196.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.263.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.50%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.114'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.41%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '530.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -14.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.759'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -13.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -13.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -12.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -13.43%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -12.42%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -15.91%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.94%  '
$ws.Range('E32').Value = '  -10.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -18.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '534.96'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -12.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.76'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -14.96%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0453'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0857'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.48%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -16.95%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.127'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -11.89%  '
$ws.Range('E42').Value = '  -18.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.933.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -11.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.267'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -13.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₃0589'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -18.37%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -15.99%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.17'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.44%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -18.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.15'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.72%  '
$ws.Range('E51').Value = '  -12.57%  '
